$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Contrasenia" column (E) for data rows 2-11: "gw" -> "silverarrow"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "silverarrow"
}

# Update the active selection to reflect the edited range, as in the source workbook
$ws.Range("E3:E11").Select()
